# Regenerate save_data to use K (strike count) instead of Strike# values.
# Only column G ("K") changes across data rows 2-49; all other columns
# (TB, PC, dS0, dSF, IP, I0, IF, date, index) remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 2
    13 = 5
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 2
    19 = 0
    20 = 3
    21 = 1
    22 = 1
    23 = 0
    24 = 2
    25 = 0
    26 = 4
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 2
    35 = 1
    36 = 1
    37 = 2
    38 = 4
    39 = 1
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 3
    45 = 1
    46 = 2
    47 = 1
    48 = 2
    49 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
